# Updated BGR model - 2025-08-15 20:01
$wb = $excel.ActiveWorkbook

# --- 1. Rename the "varbl map" sheet to "timeslice map" and repurpose it
#        as the new ~Timeslice_Map definition table ---
$wsTS = $wb.Worksheets.Item("varbl map")
$wsTS.Name = "timeslice map"

$wsTS.Range("A1").Value = "~Timeslice_Map"

$wsTS.Range("A2").Value = "dimension"
$wsTS.Range("B2").Value = "name"
$wsTS.Range("C2").Value = "description"

$wsTS.Range("A3").Value = "ts_type"
$wsTS.Range("B3").Value = "*,-s?a*"
$wsTS.Range("C3").Value = "hourly"

$wsTS.Range("A4").Value = "ts_type"
$wsTS.Range("B4").Value = "s?a*"
$wsTS.Range("C4").Value = "aggregated"

$wsTS.Range("A5").Value = "ts_season"
$wsTS.Range("B5").Value = "S1*"
$wsTS.Range("C5").Formula = "=LEFT(B5,2)"

$wsTS.Range("A6").Value = "ts_season"
$wsTS.Range("B6").Value = "S2*"
$wsTS.Range("C6").Formula = "=LEFT(B6,2)"

$wsTS.Range("A7").Value = "ts_season"
$wsTS.Range("B7").Value = "S3*"
$wsTS.Range("C7").Formula = "=LEFT(B7,2)"

$wsTS.Range("A8").Value = "ts_season"
$wsTS.Range("B8").Value = "S4*"
$wsTS.Range("C8").Formula = "=LEFT(B8,2)"

$wsTS.Range("A9").Value = "ts_season"
$wsTS.Range("B9").Value = "S5*"
$wsTS.Range("C9").Formula = "=LEFT(B9,2)"

$wsTS.Range("A10").Value = "ts_season"
$wsTS.Range("B10").Value = "S6*"
$wsTS.Range("C10").Formula = "=LEFT(B10,2)"

# --- 2. Update TS_Defs row 6 (the VAR_POUT / Power row) ---
$wsDefs = $wb.Worksheets.Item("TS_Defs")
$wsDefs.Range("C6").Value = "ELE,STG,IRE,-Grid"
$wsDefs.Range("H6").Value = "NRG"

# --- 3. Make "timeslice map" the active sheet / tab ---
$wsTS.Activate()
$wsTS.Select()
$wsTS.Range("A2").Select()

$excel.ActiveWindow.DisplayWorkbookTabs = $true
